$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1702.7142
$ws.Range("J17").Value = 1702.7142
$ws.Range("L17").Value = 5108.142599999999
$ws.Range("N17").Value = -5444.142599999999

$ws.Range("H53").Value = 403.8
$ws.Range("I53").Value = 424.6
$ws.Range("J53").Value = 383
$ws.Range("K53").Value = 424.6
$ws.Range("L53").Value = 383
$ws.Range("M53").Value = 212.4
$ws.Range("N53").Value = -1657

$ws.Range("H62").Value = 12117.091
$ws.Range("I62").Value = 1747.25
$ws.Range("K62").Value = 1747.25
$ws.Range("M62").Value = -1123.25

$ws.Range("H65").Value = 12117.091
$ws.Range("I65").Value = 1747.25
$ws.Range("K65").Value = 8736.25
$ws.Range("M65").Value = -5616.25

$ws.Range("H87").Value = 88434.5
$ws.Range("J87").Value = 88434.5
$ws.Range("L87").Value = 88434.5
$ws.Range("N87").Value = -90930.5

$ws.Range("H90").Value = 88434.5
$ws.Range("J90").Value = 88434.5
$ws.Range("L90").Value = 265303.5
$ws.Range("N90").Value = -277783.5

$ws.Range("H99").Value = 394.875
$ws.Range("I99").Value = 358.5
$ws.Range("J99").Value = 649.5
$ws.Range("K99").Value = 1075.5
$ws.Range("L99").Value = 1948.5
$ws.Range("M99").Value = 422.5
$ws.Range("N99").Value = -4944.5

$ws.Range("H113").Value = 125001340
$ws.Range("I113").Value = 50000724
$ws.Range("J113").Value = 200001970
$ws.Range("K113").Value = 50000724
$ws.Range("L113").Value = 200001970
$ws.Range("M113").Value = -49997470
$ws.Range("N113").Value = -200008478

$ws.Range("H116").Value = 5624.3335
$ws.Range("I116").Value = 4882.25
$ws.Range("K116").Value = 4882.25
$ws.Range("M116").Value = -1440.25

$ws.Range("H118").Value = 349.2857
$ws.Range("I118").Value = 349.2857
$ws.Range("K118").Value = 1047.8571
$ws.Range("M118").Value = 609.1428999999998

$ws.Range("H125").Value = 2649.077
$ws.Range("J125").Value = 3427
$ws.Range("L125").Value = 30843
$ws.Range("N125").Value = -35763

$ws.Range("H133").Value = 69963.336
$ws.Range("J133").Value = 69963.336
$ws.Range("L133").Value = 69963.336
$ws.Range("N133").Value = -80083.336

$ws.Range("H138").Value = 2198.6
$ws.Range("I138").Value = 1384
$ws.Range("J138").Value = 2417.4478
$ws.Range("K138").Value = 4152
$ws.Range("L138").Value = 7252.3434
$ws.Range("M138").Value = 988
$ws.Range("N138").Value = -17532.3434

$ws.Range("H141").Value = 13772
$ws.Range("I141").Value = 14855.375
$ws.Range("K141").Value = 44566.125
$ws.Range("M141").Value = -39386.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 463.7857
$ws.Range("I2").Value = 343.27274
$ws.Range("K2").Value = 343.27274
$ws.Range("M2").Value = -230.27274

$ws.Range("H74").Value = 9633696
$ws.Range("I74").Value = 11906629
$ws.Range("K74").Value = 11906629
$ws.Range("M74").Value = -11905755

$ws.Range("H77").Value = 9633696
$ws.Range("I77").Value = 11906629
$ws.Range("K77").Value = 59533145
$ws.Range("M77").Value = -59528777

$ws.Range("H110").Value = 1298.1904
$ws.Range("I110").Value = 1112.45
$ws.Range("K110").Value = 1112.45
$ws.Range("M110").Value = 932.55

$ws.Range("H116").Value = 463.7857
$ws.Range("I116").Value = 343.27274
$ws.Range("K116").Value = 343.27274
$ws.Range("M116").Value = 1950.72726

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 463.7857
$ws.Range("I3").Value = 343.27274
$ws.Range("K3").Value = 343.27274
$ws.Range("M3").Value = -229.27274

$ws.Range("H94").Value = 876.1429000000001
$ws.Range("I94").Value = 353.07693
$ws.Range("K94").Value = 353.07693
$ws.Range("M94").Value = 97.92307

$ws.Range("H99").Value = 3904
$ws.Range("I99").Value = 2655
$ws.Range("J99").Value = 4320.3335
$ws.Range("K99").Value = 2655
$ws.Range("L99").Value = 4320.3335
$ws.Range("M99").Value = -1157
$ws.Range("N99").Value = -7316.3335

$ws.Range("H105").Value = 2960.1667
$ws.Range("J105").Value = 2775.5334
$ws.Range("L105").Value = 2775.5334
$ws.Range("N105").Value = -6269.5334

$ws.Range("H107").Value = 3863.7144
$ws.Range("I107").Value = 4786.8
$ws.Range("J107").Value = 1556
$ws.Range("K107").Value = 4786.8
$ws.Range("L107").Value = 1556
$ws.Range("M107").Value = -2866.8
$ws.Range("N107").Value = -5396

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1347.2413
$ws.Range("I107").Value = 435.4737
$ws.Range("J107").Value = 3079.6
$ws.Range("K107").Value = 435.4737
$ws.Range("L107").Value = 3079.6
$ws.Range("M107").Value = 1484.5263
$ws.Range("N107").Value = -6919.6

$ws.Range("H122").Value = 3611.5
$ws.Range("I122").Value = 3611.5
$ws.Range("K122").Value = 10834.5
$ws.Range("M122").Value = -8384.5

$ws.Range("H132").Value = 4999.4
$ws.Range("I132").Value = 4249.75
$ws.Range("K132").Value = 12749.25
$ws.Range("M132").Value = -10219.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 267.8125
$ws.Range("I2").Value = 266.14285
$ws.Range("K2").Value = 1596.8571
$ws.Range("M2").Value = -1483.8571

$ws.Range("H12").Value = 1393.4
$ws.Range("J12").Value = 540
$ws.Range("L12").Value = 1620
$ws.Range("N12").Value = -1966

$ws.Range("H109").Value = 1906.2
$ws.Range("I109").Value = 1906.2
$ws.Range("K109").Value = 5718.6
$ws.Range("M109").Value = -4678.6

$ws.Range("H114").Value = 780.5
$ws.Range("J114").Value = 399.66666
$ws.Range("L114").Value = 1198.99998
$ws.Range("N114").Value = -7706.999980000001

$ws.Range("H117").Value = 1099
$ws.Range("J117").Value = 1099
$ws.Range("L117").Value = 3297
$ws.Range("N117").Value = -10181

$ws.Range("H136").Value = 4770.25
$ws.Range("I136").Value = 4770.25
$ws.Range("K136").Value = 14310.75
$ws.Range("M136").Value = -9210.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2386.64
$ws.Range("I102").Value = 1686.0588
$ws.Range("K102").Value = 1686.0588
$ws.Range("M102").Value = -64.05880000000002

$ws.Range("H113").Value = 3657.1538
$ws.Range("I113").Value = 2768.4
$ws.Range("J113").Value = 4212.625
$ws.Range("K113").Value = 2768.4
$ws.Range("L113").Value = 4212.625
$ws.Range("M113").Value = -598.4000000000001
$ws.Range("N113").Value = -8552.625

$ws.Range("H122").Value = 2074.5454
$ws.Range("I122").Value = 1625.625
$ws.Range("K122").Value = 4876.875
$ws.Range("M122").Value = -2426.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2329.3
$ws.Range("I22").Value = 2857.8
$ws.Range("J22").Value = 1800.8
$ws.Range("K22").Value = 2857.8
$ws.Range("L22").Value = 1800.8
$ws.Range("M22").Value = -2562.8
$ws.Range("N22").Value = -2390.8

$ws.Range("H27").Value = 2329.3
$ws.Range("I27").Value = 2857.8
$ws.Range("J27").Value = 1800.8
$ws.Range("K27").Value = 2857.8
$ws.Range("L27").Value = 1800.8
$ws.Range("M27").Value = -2750.8
$ws.Range("N27").Value = -2014.8

$ws.Range("H46").Value = 3967.2424
$ws.Range("I46").Value = 1977.6154
$ws.Range("J46").Value = 11357.286
$ws.Range("K46").Value = 1977.6154
$ws.Range("L46").Value = 11357.286
$ws.Range("M46").Value = -1789.6154
$ws.Range("N46").Value = -11733.286

$ws.Range("H55").Value = 32258668
$ws.Range("J55").Value = 684.5
$ws.Range("L55").Value = 684.5
$ws.Range("N55").Value = -1030.5

$ws.Range("H122").Value = 5854.074
$ws.Range("I122").Value = 5293.125
$ws.Range("K122").Value = 15879.375
$ws.Range("M122").Value = -13429.375

$ws.Range("H136").Value = 67096.52
$ws.Range("I136").Value = 8643.714
$ws.Range("K136").Value = 25931.142
$ws.Range("M136").Value = -23381.142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 10000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 10000
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -9707
$ws.Range("N26").ClearContents()

$ws.Range("H132").Value = 1436717
$ws.Range("I132").Value = 9403
$ws.Range("J132").Value = 5005002
$ws.Range("K132").Value = 28209
$ws.Range("L132").Value = 15015006
$ws.Range("M132").Value = -25679
$ws.Range("N132").Value = -15020066
